# Apply the crypto price/volume update for Sun Feb 26 05:30:20 UTC 2023 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Update per-row Price (D) and Volume(1h) (E) values
Set-TextValue "D2" "23.189.67"
Set-TextValue "E2" "  +0.42%  "

Set-TextValue "D3" "1.601.82"
Set-TextValue "E3" "  +0.18%  "

Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.13%  "

Set-TextValue "D5" "1.000"
Set-TextValue "E5" "  -0.11%  "

Set-TextValue "D6" "303.04"
Set-TextValue "E6" "  +0.70%  "

Set-TextValue "D7" "0.3780"
Set-TextValue "E7" "  +0.18%  "

Set-TextValue "D8" "51.77"
Set-TextValue "E8" "  +3.20%  "

Set-TextValue "D9" "0.3613"
Set-TextValue "E9" "  -0.66%  "

Set-TextValue "D10" "1.264"
Set-TextValue "E10" "  +0.87%  "

Set-TextValue "E11" "  -0.15%  "

Set-TextValue "D12" "0.08123"
Set-TextValue "E12" "  +0.03%  "

Set-TextValue "D13" "22.59"
Set-TextValue "E13" "  -1.58%  "

Set-TextValue "D14" "6.582"
Set-TextValue "E14" "  +0.35%  "

Set-TextValue "D15" "7.390"
Set-TextValue "E15" "  +0.63%  "

Set-TextValue "D16" "0.00001248"
Set-TextValue "E16" "  -0.50%  "

Set-TextValue "D17" "1.601.41"
Set-TextValue "E17" "  +0.03%  "

Set-TextValue "D18" "93.79"
Set-TextValue "E18" "  +2.43%  "

Set-TextValue "D19" "0.06883"
Set-TextValue "E19" "  +0.13%  "

Set-TextValue "D20" "18.04"
Set-TextValue "E20" "  -0.94%  "

Set-TextValue "D21" "6.533"
Set-TextValue "E21" "  -0.18%  "

Set-TextValue "E22" "  -0.03%  "

Set-TextValue "D23" "12.95"
Set-TextValue "E23" "  -0.12%  "

Set-TextValue "D24" "23.196.51"
Set-TextValue "E24" "  +0.44%  "

Set-TextValue "D25" "2.387"
Set-TextValue "E25" "  +1.93%  "

Set-TextValue "D26" "2.978"
Set-TextValue "E26" "  +9.81%  "

Set-TextValue "D27" "21.18"
Set-TextValue "E27" "  +0.70%  "

Set-TextValue "D28" "149.55"
Set-TextValue "E28" "  -0.25%  "

Set-TextValue "D29" "5.245"
Set-TextValue "E29" "  -0.08%  "

Set-TextValue "D30" "133.81"
Set-TextValue "E30" "  +1.58%  "

Set-TextValue "D31" "2.385"
Set-TextValue "E31" "  -1.85%  "

Set-TextValue "D32" "6.849"
Set-TextValue "E32" "  +0.73%  "

Set-TextValue "D33" "1.780.31"
Set-TextValue "E33" "  +0.17%  "

Set-TextValue "D34" "0.9726"
Set-TextValue "E34" "  +2.58%  "

Set-TextValue "D35" "0.07520"
Set-TextValue "E35" "  -1.68%  "

Set-TextValue "D36" "10.31"
Set-TextValue "E36" "  +2.83%  "

Set-TextValue "D37" "0.02717"

Set-TextValue "D40" "0.08800"
Set-TextValue "E40" "  -1.11%  "

Set-TextValue "D41" "0.7097"
Set-TextValue "E41" "  +0.37%  "

Set-TextValue "E42" "  -0.62%  "

Set-TextValue "D43" "12.50"
Set-TextValue "E43" "  -0.90%  "

Set-TextValue "D44" "15.39"
Set-TextValue "E44" "  -0.49%  "

Set-TextValue "D45" "0.6538"
Set-TextValue "E45" "  -0.84%  "

Set-TextValue "D46" "2.307"
Set-TextValue "E46" "  +0.69%  "

Set-TextValue "E47" "  +0.85%  "

Set-TextValue "D48" "132.31"
Set-TextValue "E48" "  +0.59%  "

Set-TextValue "D49" "0.07958"
Set-TextValue "E49" "  +0.13%  "

Set-TextValue "E50" "  +0.06%  "

Set-TextValue "D51" "1.222"
Set-TextValue "E51" "  +3.48%  "

# Rows 38 and 39 swap rankings: row 38 becomes Algorand, row 39 becomes InternetComputer(DFINITY)
Set-TextValue "B38" "Algorand"
Set-TextValue "C38" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D38" "0.2508"
Set-TextValue "E38" "  -1.34%  "

Set-TextValue "B39" "InternetComputer(DFINITY)"
Set-TextValue "C39" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D39" "6.117"
Set-TextValue "E39" "  -1.55%  "
